$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 117 / 118: swap the full data payload (B..AD) between the two rows.
# Row index column A (and the row structure itself) stays put; only the
# match data moved from one row to the other.
$row117 = $ws.Range("B117:AD117").Value2
$row118 = $ws.Range("B118:AD118").Value2

$ws.Range("B117:AD117").Value2 = $row118
$ws.Range("B118:AD118").Value2 = $row117

# --- Rows 234 / 235 / 237: 3-way cyclic rotation of the data payload
# (234 <- 237, 235 <- 234, 237 <- 235), keeping column A fixed. Columns I/J
# are unused on these rows, so only B..H and K..AD are touched.
$row234_bh = $ws.Range("B234:H234").Value2
$row234_kz = $ws.Range("K234:AD234").Value2
$row235_bh = $ws.Range("B235:H235").Value2
$row235_kz = $ws.Range("K235:AD235").Value2
$row237_bh = $ws.Range("B237:H237").Value2
$row237_kz = $ws.Range("K237:AD237").Value2

$ws.Range("B234:H234").Value2 = $row237_bh
$ws.Range("K234:AD234").Value2 = $row237_kz

$ws.Range("B235:H235").Value2 = $row234_bh
$ws.Range("K235:AD235").Value2 = $row234_kz

$ws.Range("B237:H237").Value2 = $row235_bh
$ws.Range("K237:AD237").Value2 = $row235_kz
